# Generate Report for Handoff
#
# Re-run of the localization-status report generator: the row previously
# reporting "e85bc520-..." as freshly handed off has been superseded by a
# new handoff pass. "2d5489b1-..." and "7de476d6-..." (already "Ready for
# handoff") move up to rows 7-8, and "e85bc520-..." drops to row 9 carrying
# a brand-new handoff timestamp and a status of "Ready for handoff".
#
# Hyperlink *targets* (Address / r:id) stay anchored to their row position;
# only the cell text and the hyperlink *display* text move.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @{ Row = 7; A = "2d5489b1-3d79-483c-a379-dd25fef59254.md" },
    @{ Row = 8; A = "7de476d6-9618-4066-b0f5-43ac4d0276d1.md" },
    @{ Row = 9; A = "e85bc520-d944-4418-b933-05784aa62d3f.md" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = "Ready for handoff"
    $ws.Range("C$row").Value = "Ready for handoff"
}

$overviewLinks = @{
    '$A$7' = "2d5489b1-3d79-483c-a379-dd25fef59254.md";
    '$A$8' = "7de476d6-9618-4066-b0f5-43ac4d0276d1.md";
    '$A$9' = "e85bc520-d944-4418-b933-05784aa62d3f.md"
}

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($overviewLinks.ContainsKey($addr)) {
        $h.TextToDisplay = $overviewLinks[$addr]
    }
}

# ---- Per-locale sheets (zh-cn, de-de) -------------------------------------
$locales = @(
    @{
        Sheet = "zh-cn";
        Rows = @(
            @{ Row = 7; A = "2d5489b1-3d79-483c-a379-dd25fef59254.md"; C = "2d5489b1-3d79-483c-a379-dd25fef59254.9b49fb7f1246cb224573d0050187ee3e6945911c.zh-cn.xlf"; D = "2016-03-09 09:48:47" },
            @{ Row = 8; A = "7de476d6-9618-4066-b0f5-43ac4d0276d1.md"; C = "7de476d6-9618-4066-b0f5-43ac4d0276d1.b4c6246186f9eca2e193bcd5a6dd31e66525e085.zh-cn.xlf"; D = "2016-03-09 09:46:56" },
            @{ Row = 9; A = "e85bc520-d944-4418-b933-05784aa62d3f.md"; C = "e85bc520-d944-4418-b933-05784aa62d3f.154992c9893038fe11bb830932967dfb491ab6d5.zh-cn.xlf"; D = "2016-03-09 09:51:11" }
        )
    },
    @{
        Sheet = "de-de";
        Rows = @(
            @{ Row = 7; A = "2d5489b1-3d79-483c-a379-dd25fef59254.md"; C = "2d5489b1-3d79-483c-a379-dd25fef59254.9b49fb7f1246cb224573d0050187ee3e6945911c.de-de.xlf"; D = "2016-03-09 09:48:50" },
            @{ Row = 8; A = "7de476d6-9618-4066-b0f5-43ac4d0276d1.md"; C = "7de476d6-9618-4066-b0f5-43ac4d0276d1.b4c6246186f9eca2e193bcd5a6dd31e66525e085.de-de.xlf"; D = "2016-03-09 09:46:59" },
            @{ Row = 9; A = "e85bc520-d944-4418-b933-05784aa62d3f.md"; C = "e85bc520-d944-4418-b933-05784aa62d3f.154992c9893038fe11bb830932967dfb491ab6d5.de-de.xlf"; D = "2016-03-09 09:51:15" }
        )
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    $linkMap = @{}
    foreach ($r in $locale.Rows) {
        $row = $r.Row
        $ws.Range("A$row").Value = $r.A
        $ws.Range("B$row").Value = "Ready for handoff"
        $ws.Range("C$row").Value = $r.C
        $ws.Range("D$row").Value = $r.D

        $linkMap["`$A`$$row"] = $r.A
        $linkMap["`$C`$$row"] = $r.C
    }

    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($linkMap.ContainsKey($addr)) {
            $h.TextToDisplay = $linkMap[$addr]
        }
    }
}
